$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the <cols> definition: column 1 should only cover column 1 (the
#    original file had an overlapping "min=1 max=2" range that is redundant
#    with the following "min=2 max=2" entry). Re-setting the width on
#    column 1 alone causes the overlap to be cleaned up on save.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.7109375

# ---------------------------------------------------------------------------
# 2. Insert two new rows right after row 12 ("Docentes responsaveis:"),
#    pushing the old rows 13-21 down to rows 15-23. This also keeps every
#    existing row's custom height attached to the correct (shifted) row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).EntireRow.Insert()
$ws.Rows.Item(13).EntireRow.Insert()

# The inserted rows copied the formatting of the row above (row 12, which
# only has a bold "A" column style). Clear the stray A-column formatting and
# paint the B/C columns with the normal/red formatting used throughout the
# sheet by copying the format from an already well-formatted row.
$ws.Cells.Item(13,1).Clear()
$ws.Cells.Item(14,1).Clear()

$ws.Range("B18:C18").Copy() | Out-Null
$ws.Range("B13:C14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Remove the explicit row height picked up from the row-insert operation so
# the two new rows use the sheet's default height, matching the target.
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()

# ---------------------------------------------------------------------------
# 3. Fill in the two new "Docentes responsaveis" rows (13 and 14).
# ---------------------------------------------------------------------------
$ws.Cells.Item(13,2).Value = '5817344 - Livia Melo Carneiro'
$ws.Cells.Item(13,3).Value = '5817344 - Livia Melo Carneiro'

$ws.Cells.Item(14,2).Value = '6310296 - Patrícia Caroline Molgero Da Rós'
$ws.Cells.Item(14,3).Value = '6310296 - Patrícia Caroline Molgero Da Rós'

# ---------------------------------------------------------------------------
# 4. Row 10 ("Objetivos:") previously held the wrong text (a docente name).
#    Replace it with the actual Portuguese objectives text.
# ---------------------------------------------------------------------------
$objetivos = 'Dar o embasamento dos conceitos elementares em química aos alunos, capacitando-os para o prosseguimento dos estudos nas disciplinas correlatas posteriores, principalmente quanto aos conceitos da estrutura atômica; das ligações química e forças intermoleculares; da geometria das moléculas; da natureza dos compostos; das reações químicas em solução aquosa, tanto de dupla-troca como de oxirredução; das propriedades do estado gasoso e das soluções e da estequiometria e cálculos em química, com ênfase em casos contendo reagentes limitantes, pureza de reagentes e rendimento de reação.'
$ws.Cells.Item(10,2).Value = $objetivos
$ws.Cells.Item(10,3).Value = $objetivos

# ---------------------------------------------------------------------------
# 5. Row 15 ("Programa resumido:", was old row 13) previously held the wrong
#    text (a date). Replace it with the actual Portuguese short syllabus
#    text.
# ---------------------------------------------------------------------------
$programaResumido = 'Sistemas de Unidades. Estrutura Atômica. Tabela Periódica. Ligação Química. Nomenclatura de compostos inorgânicos.  Definições de ácidos e bases. Forças intermoleculares.  Soluções. Gases. Reações químicas em solução aquosa. Estequiometria e Cálculos em Química.'
$ws.Cells.Item(15,2).Value = $programaResumido
$ws.Cells.Item(15,3).Value = $programaResumido

# ---------------------------------------------------------------------------
# 6. Row 17 ("Programa:", was old row 15) previously held the wrong text (a
#    docente name). Replace it with the actual Portuguese full syllabus
#    text.
# ---------------------------------------------------------------------------
$programa = 'Sistemas de unidades: Definição das Unidades mais usadas em Engenharia e transformações entre sistemas. Estrutura atômica: Natureza elétrica da matéria. A carga do elétron. O núcleo do átomo. Teoria quântica: A radiação, os quanta e os fótons. Espectros de emissão e de absorção atômica. A dualidade onda-partícula da matéria. O princípio da incerteza. Os orbitais atômicos. Os números quânticos. Configuração eletrônica dos elementos. Partículas Elementares. Tabela periódica: A Lei e a tabela Periódica. Propriedades periódicas dos elementos, átomos e íons. Ligação Química: A ligação covalente. Estrutura de Lewis. Orbitais moleculares: Limitações da teoria de ligação de valência. Hibridização. Polaridade da ligação. Geometria molecular (Modelo VSEPR). Ligação Iônica. A classificação dos sólidos. As propriedades das ligações. Os compostos de coordenação. Complexos metálicos (teoria do campo cristalino). Ligação Metálica.Nomenclatura de compostos inorgânicos: Funções Inorgânicas: ácidos; bases; sais; óxidos e nomenclaturas.Definições de ácidos e bases: Ácidos e bases (Arrhenius, Bronsted-Lowry e Lewis). Forças intermoleculares: Forças intermoleculares, líquidos e sólidosSoluções: Natureza das soluções. Dispersões coloidais e suspensões. Propriedades físicas e químicas. Tipos de soluções. Unidades e cálculos de concentração (Molaridade, fração molar, ppm, normalidade, molalidade, diluição). O processo de dissolução. Calor de dissolução. Solubilidade e temperatura.Gases (ideais e reais): Variáveis de estado. Lei combinada dos gases. Experiência de Torriceli. Pressão parcial dos gases. Teoria cinética dos gases. Gás ideal e real. Princípio de Avogadro. Reações químicas em solução aquosa: Principais reações químicas (ácido-base, precipitação, óxido-redução e complexação). Exemplos de reações formadoras de gases. Princípios de titulações ácido-base e de óxido-redução.Estequiometria e cálculos em química: Balanceamento de reações, cálculos estequiométricos, reagentes limitantes e rendimentos.'
$ws.Cells.Item(17,2).Value = $programa
$ws.Cells.Item(17,3).Value = $programa

# ---------------------------------------------------------------------------
# 7. Row 20 ("Método:", was old row 18) previously held the wrong text (a
#    docente name). Replace it with the actual method text.
# ---------------------------------------------------------------------------
$ws.Cells.Item(20,2).Value = 'Duas provas escritas'
$ws.Cells.Item(20,3).Value = 'Duas provas escritas'

# ---------------------------------------------------------------------------
# 8. Row 21 ("Critério:", was old row 19) previously held the wrong text
#    (the "Método" text). Replace it with the real evaluation criteria
#    text.
# ---------------------------------------------------------------------------
$criterio = 'A média para a primeira avaliação será calculada a partir das notas das duas provas, P1 e P2, segundo a fórmula: M1=(P1+2xP2)/3. Alunos com nota final igual ou superior a 5,0 estão aprovados; inferior a 5,0 e igual ou superior a 3,0 estão de recuperação;'
$ws.Cells.Item(21,2).Value = $criterio
$ws.Cells.Item(21,3).Value = $criterio

# ---------------------------------------------------------------------------
# 9. Row 22 ("Norma de recuperação:", was old row 20) previously held the
#    wrong text (the "Critério" text). Replace it with the real recovery
#    norm text.
# ---------------------------------------------------------------------------
$norma = 'A recuperação consistirá de uma prova envolvendo o assunto do semestre todo, à qual será atribuída nota NR. A média da segunda avaliação será calculada segunda a fórmula: M2=(M1+NR)/2. Alunos com nota M2 igual ou superior a 5,0 estarão aprovados, inferior a 5,0 estarão reprovados.'
$ws.Cells.Item(22,2).Value = $norma
$ws.Cells.Item(22,3).Value = $norma

# ---------------------------------------------------------------------------
# 10. Row 23 ("Bibliografia:", was old row 21) previously held the wrong
#     text (the "Norma de recuperação" text). Replace it with the real
#     bibliography text.
# ---------------------------------------------------------------------------
$bibliografia = 'ATKINS, Peter., Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.RUSSEL, J.B. Química geral. São Paulo: MacGrall-Hill'
$ws.Cells.Item(23,2).Value = $bibliografia
$ws.Cells.Item(23,3).Value = $bibliografia

# ---------------------------------------------------------------------------
# The worksheet dimension (A1:C21 -> A1:C23) is recalculated automatically
# by Excel when the file is saved, so nothing further is required.
# ---------------------------------------------------------------------------
